$d = $word.ActiveDocument

$replacements = @(
    @("2025-12-15 Monday", "2025-12-16 Tuesday"),
    @("72×35=", "16×72="),
    @("91×85=", "99×18="),
    @("58×90=", "43×95="),
    @("48×27=", "13×83="),
    @("79×50=", "58×74="),
    @("34×35=", "44×22="),
    @("60×42=", "65×82="),
    @("42×15=", "54×51="),
    @("47×47=", "96×51="),
    @("98×31=", "72×68="),
    @("60×56=", "63×49="),
    @("31×55=", "66×92="),
    @("39×81=", "30×61="),
    @("24×12=", "43×76="),
    @("62×52=", "91×66="),
    @("16×50=", "62×26="),
    @("77×25=", "56×46="),
    @("90×90=", "94×43="),
    @("88×46=", "30×76="),
    @("15×74=", "58×95="),
    @("86×29=", "14×61="),
    @("32×49=", "87×88="),
    @("53×18=", "66×87="),
    @("42×67=", "47×67="),
    @("85×64=", "22×21=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
